# Update gh-pages to output generated at 456a3b4
# Increments the "想去人数" (F column) counts for several events across
# the 展览, 演出 and 全部类型 sheets.

$wb = $excel.ActiveWorkbook

# Sheet: 展览 (F column values)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F11").Value = 479
$ws1.Range("F13").Value = 1119
$ws1.Range("F15").Value = 4644
$ws1.Range("F20").Value = 3638
$ws1.Range("F30").Value = 193
$ws1.Range("F32").Value = 85
$ws1.Range("F36").Value = 5939
$ws1.Range("F38").Value = 446
$ws1.Range("F46").Value = 2086

# Sheet: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F23").Value = 775

# Sheet: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F15").Value = 479
$ws4.Range("F16").Value = 1119
$ws4.Range("F18").Value = 4644
$ws4.Range("F23").Value = 3638
$ws4.Range("F28").Value = 193
$ws4.Range("F30").Value = 85
$ws4.Range("F35").Value = 5939
$ws4.Range("F37").Value = 446
$ws4.Range("F46").Value = 2086
